$d = $word.ActiveDocument

# Paragraph 4 currently holds "是小慧慧的生日哟。" and the _GoBack bookmark.
$p4 = $d.Paragraphs(4)

# 1) Its own paragraph-mark run fonts hint flips from "default" to "eastAsia".
$p4.Range.Font.NameFarEast = $p4.Range.Font.NameFarEast

# Insert two new paragraphs after it (after the bookmark, at the very end
# of the paragraph's range) for the new diary entry.
$r = $p4.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.Text = "2020年11月5日 天气晴朗"

$r5 = $p5.Range
$r5.Collapse(0)
$r5.InsertParagraphAfter()

$p6 = $d.Paragraphs(6)
$p6.Range.Text = "是我们刚开始在一起的第一天。"

Write-Output $d.Paragraphs.Count
